# Update sheet1 (LP1912) header info
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value2 = "Última actualización: 10:50:41"
$ws1.Range("A3").Value2 = "Total filas: 109"

# Update / insert data rows 87-114 on sheet1 with the new scrape results
$ws1.Cells.Item(87,1).Value2 = "10:50:41"
$ws1.Cells.Item(87,2).Value2 = "10:59"
$ws1.Cells.Item(87,3).Value2 = "10_OLMOS"
$ws1.Cells.Item(87,4).Value2 = 9
$ws1.Cells.Item(87,5).Value2 = "LP1912"
$ws1.Cells.Item(88,1).Value2 = "09:23:23"
$ws1.Cells.Item(88,2).Value2 = "11:01"
$ws1.Cells.Item(88,3).Value2 = "81_EL PELIGRO"
$ws1.Cells.Item(88,4).Value2 = 98
$ws1.Cells.Item(88,5).Value2 = "LP1912"
$ws1.Cells.Item(89,1).Value2 = "10:05:51"
$ws1.Cells.Item(89,2).Value2 = "11:04"
$ws1.Cells.Item(89,3).Value2 = "23_HERNANDEZ"
$ws1.Cells.Item(89,4).Value2 = 59
$ws1.Cells.Item(89,5).Value2 = "LP1912"
$ws1.Cells.Item(90,1).Value2 = "10:37:52"
$ws1.Cells.Item(90,2).Value2 = "11:06"
$ws1.Cells.Item(90,3).Value2 = "23_HERNANDEZ"
$ws1.Cells.Item(90,4).Value2 = 29
$ws1.Cells.Item(90,5).Value2 = "LP1912"
$ws1.Cells.Item(91,1).Value2 = "09:23:23"
$ws1.Cells.Item(91,2).Value2 = "11:10"
$ws1.Cells.Item(91,3).Value2 = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(91,4).Value2 = 107
$ws1.Cells.Item(91,5).Value2 = "LP1912"
$ws1.Cells.Item(92,1).Value2 = "09:23:23"
$ws1.Cells.Item(92,2).Value2 = "11:14"
$ws1.Cells.Item(92,3).Value2 = "14_ABASTO"
$ws1.Cells.Item(92,4).Value2 = 111
$ws1.Cells.Item(92,5).Value2 = "LP1912"
$ws1.Cells.Item(93,1).Value2 = "09:23:23"
$ws1.Cells.Item(93,2).Value2 = "11:15"
$ws1.Cells.Item(93,3).Value2 = "15X38_ABASTO"
$ws1.Cells.Item(93,4).Value2 = 112
$ws1.Cells.Item(93,5).Value2 = "LP1912"
$ws1.Cells.Item(94,1).Value2 = "10:37:52"
$ws1.Cells.Item(94,2).Value2 = "11:25"
$ws1.Cells.Item(94,3).Value2 = "16_SANTA ANA"
$ws1.Cells.Item(94,4).Value2 = 48
$ws1.Cells.Item(94,5).Value2 = "LP1912"
$ws1.Cells.Item(95,1).Value2 = "10:05:51"
$ws1.Cells.Item(95,2).Value2 = "11:28"
$ws1.Cells.Item(95,3).Value2 = "10_OLMOS"
$ws1.Cells.Item(95,4).Value2 = 83
$ws1.Cells.Item(95,5).Value2 = "LP1912"
$ws1.Cells.Item(96,1).Value2 = "10:50:41"
$ws1.Cells.Item(96,2).Value2 = "11:29"
$ws1.Cells.Item(96,3).Value2 = "10_OLMOS"
$ws1.Cells.Item(96,4).Value2 = 39
$ws1.Cells.Item(96,5).Value2 = "LP1912"
$ws1.Cells.Item(97,1).Value2 = "10:05:51"
$ws1.Cells.Item(97,2).Value2 = "11:30"
$ws1.Cells.Item(97,3).Value2 = "215C_EL PATO"
$ws1.Cells.Item(97,4).Value2 = 85
$ws1.Cells.Item(97,5).Value2 = "LP1912"
$ws1.Cells.Item(98,1).Value2 = "10:05:51"
$ws1.Cells.Item(98,2).Value2 = "11:31"
$ws1.Cells.Item(98,3).Value2 = "16_SANTA ANA"
$ws1.Cells.Item(98,4).Value2 = 86
$ws1.Cells.Item(98,5).Value2 = "LP1912"
$ws1.Cells.Item(99,1).Value2 = "10:05:51"
$ws1.Cells.Item(99,2).Value2 = "11:41"
$ws1.Cells.Item(99,3).Value2 = "215B_EL PATO"
$ws1.Cells.Item(99,4).Value2 = 96
$ws1.Cells.Item(99,5).Value2 = "LP1912"
$ws1.Cells.Item(100,1).Value2 = "10:05:51"
$ws1.Cells.Item(100,2).Value2 = "11:45"
$ws1.Cells.Item(100,3).Value2 = "15X38_ABASTO"
$ws1.Cells.Item(100,4).Value2 = 100
$ws1.Cells.Item(100,5).Value2 = "LP1912"
$ws1.Cells.Item(101,1).Value2 = "10:05:51"
$ws1.Cells.Item(101,2).Value2 = "11:52"
$ws1.Cells.Item(101,3).Value2 = "225_GOMEZ"
$ws1.Cells.Item(101,4).Value2 = 107
$ws1.Cells.Item(101,5).Value2 = "LP1912"
$ws1.Cells.Item(102,1).Value2 = "10:37:52"
$ws1.Cells.Item(102,2).Value2 = "11:53"
$ws1.Cells.Item(102,3).Value2 = "23_HERNANDEZ"
$ws1.Cells.Item(102,4).Value2 = 76
$ws1.Cells.Item(102,5).Value2 = "LP1912"
$ws1.Cells.Item(103,1).Value2 = "10:50:41"
$ws1.Cells.Item(103,2).Value2 = "11:53"
$ws1.Cells.Item(103,3).Value2 = "225_GOMEZ"
$ws1.Cells.Item(103,4).Value2 = 63
$ws1.Cells.Item(103,5).Value2 = "LP1912"
$ws1.Cells.Item(104,1).Value2 = "10:50:41"
$ws1.Cells.Item(104,2).Value2 = "11:54"
$ws1.Cells.Item(104,3).Value2 = "23_HERNANDEZ"
$ws1.Cells.Item(104,4).Value2 = 64
$ws1.Cells.Item(104,5).Value2 = "LP1912"
$ws1.Cells.Item(105,1).Value2 = "10:05:51"
$ws1.Cells.Item(105,2).Value2 = "11:58"
$ws1.Cells.Item(105,3).Value2 = "17_ROMERO"
$ws1.Cells.Item(105,4).Value2 = 113
$ws1.Cells.Item(105,5).Value2 = "LP1912"
$ws1.Cells.Item(106,1).Value2 = "10:37:52"
$ws1.Cells.Item(106,2).Value2 = "12:05"
$ws1.Cells.Item(106,3).Value2 = "11_ETCHEVERRY"
$ws1.Cells.Item(106,4).Value2 = 88
$ws1.Cells.Item(106,5).Value2 = "LP1912"
$ws1.Cells.Item(107,1).Value2 = "10:37:52"
$ws1.Cells.Item(107,2).Value2 = "12:10"
$ws1.Cells.Item(107,3).Value2 = "15_ABASTO"
$ws1.Cells.Item(107,4).Value2 = 93
$ws1.Cells.Item(107,5).Value2 = "LP1912"
$ws1.Cells.Item(108,1).Value2 = "10:37:52"
$ws1.Cells.Item(108,2).Value2 = "12:10"
$ws1.Cells.Item(108,3).Value2 = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(108,4).Value2 = 93
$ws1.Cells.Item(108,5).Value2 = "LP1912"
$ws1.Cells.Item(109,1).Value2 = "10:37:52"
$ws1.Cells.Item(109,2).Value2 = "12:16"
$ws1.Cells.Item(109,3).Value2 = "10_OLMOS"
$ws1.Cells.Item(109,4).Value2 = 99
$ws1.Cells.Item(109,5).Value2 = "LP1912"
$ws1.Cells.Item(110,1).Value2 = "10:37:52"
$ws1.Cells.Item(110,2).Value2 = "12:21"
$ws1.Cells.Item(110,3).Value2 = "215C_EL PATO"
$ws1.Cells.Item(110,4).Value2 = 104
$ws1.Cells.Item(110,5).Value2 = "LP1912"
$ws1.Cells.Item(111,1).Value2 = "10:37:52"
$ws1.Cells.Item(111,2).Value2 = "12:32"
$ws1.Cells.Item(111,3).Value2 = "14_ABASTO"
$ws1.Cells.Item(111,4).Value2 = 115
$ws1.Cells.Item(111,5).Value2 = "LP1912"
$ws1.Cells.Item(112,1).Value2 = "10:37:52"
$ws1.Cells.Item(112,2).Value2 = "12:34"
$ws1.Cells.Item(112,3).Value2 = "15_ABASTO"
$ws1.Cells.Item(112,4).Value2 = 117
$ws1.Cells.Item(112,5).Value2 = "LP1912"
$ws1.Cells.Item(113,1).Value2 = "10:50:41"
$ws1.Cells.Item(113,2).Value2 = "12:36"
$ws1.Cells.Item(113,3).Value2 = "27_EL RETIRO"
$ws1.Cells.Item(113,4).Value2 = 106
$ws1.Cells.Item(113,5).Value2 = "LP1912"
$ws1.Cells.Item(114,1).Value2 = "10:50:41"
$ws1.Cells.Item(114,2).Value2 = "12:48"
$ws1.Cells.Item(114,3).Value2 = "16_SANTA ANA"
$ws1.Cells.Item(114,4).Value2 = 118
$ws1.Cells.Item(114,5).Value2 = "LP1912"

# Update "Ultima actualizacion" timestamp on the other two sheets
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value2 = "Última actualización: 10:50:41"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value2 = "Última actualización: 10:50:41"
